# Adds 13 new entries (row 817-829) dated 46063 (2026-02-10) to the Wellness log,
# mirroring the existing table layout, and extends the I-column Volume*Intensite formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 817
$ws.Range("A814:H814").Copy() | Out-Null
$ws.Range("A817:H817").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(817,1).Value = 46063
$ws.Cells.Item(817,2).Value = "Nathanael Beta"
$ws.Cells.Item(817,3).Value = 70
$ws.Cells.Item(817,4).Value = 7
$ws.Cells.Item(817,5).Value = 6
$ws.Cells.Item(817,6).Value = 5
$ws.Cells.Item(817,7).Value = "Dos"
$ws.Cells.Item(817,8).Value = 6
$ws.Cells.Item(817,9).Formula = "=C817*D817"

# Row 818
$ws.Range("A813:H813").Copy() | Out-Null
$ws.Range("A818:H818").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(818,1).Value = 46063
$ws.Cells.Item(818,2).Value = "Omar Benyounes"
$ws.Cells.Item(818,3).Value = 70
$ws.Cells.Item(818,4).Value = 8
$ws.Cells.Item(818,5).Value = 6
$ws.Cells.Item(818,6).Value = 0
$ws.Cells.Item(818,8).Value = 8
$ws.Cells.Item(818,9).Formula = "=C818*D818"

# Row 819
$ws.Range("A813:H813").Copy() | Out-Null
$ws.Range("A819:H819").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(819,1).Value = 46063
$ws.Cells.Item(819,2).Value = "Theo Owono"
$ws.Cells.Item(819,3).Value = 70
$ws.Cells.Item(819,4).Value = 8
$ws.Cells.Item(819,5).Value = 6
$ws.Cells.Item(819,6).Value = 0
$ws.Cells.Item(819,8).Value = 4
$ws.Cells.Item(819,9).Formula = "=C819*D819"

# Row 820
$ws.Range("A814:H814").Copy() | Out-Null
$ws.Range("A820:H820").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(820,1).Value = 46063
$ws.Cells.Item(820,2).Value = "Yoann Martelat"
$ws.Cells.Item(820,3).Value = 70
$ws.Cells.Item(820,4).Value = 9
$ws.Cells.Item(820,5).Value = 8
$ws.Cells.Item(820,6).Value = 5
$ws.Cells.Item(820,7).Value = "Genou"
$ws.Cells.Item(820,8).Value = 7
$ws.Cells.Item(820,9).Formula = "=C820*D820"

# Row 821
$ws.Range("A814:H814").Copy() | Out-Null
$ws.Range("A821:H821").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(821,1).Value = 46063
$ws.Cells.Item(821,2).Value = "Kamal Bafounta"
$ws.Cells.Item(821,3).Value = 70
$ws.Cells.Item(821,4).Value = 8
$ws.Cells.Item(821,5).Value = 5
$ws.Cells.Item(821,6).Value = 3
$ws.Cells.Item(821,7).Value = "Genou"
$ws.Cells.Item(821,8).Value = 7
$ws.Cells.Item(821,9).Formula = "=C821*D821"

# Row 822
$ws.Range("A813:H813").Copy() | Out-Null
$ws.Range("A822:H822").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(822,1).Value = 46063
$ws.Cells.Item(822,2).Value = "Naim Ighbane"
$ws.Cells.Item(822,3).Value = 70
$ws.Cells.Item(822,4).Value = 9
$ws.Cells.Item(822,5).Value = 7
$ws.Cells.Item(822,6).Value = 0
$ws.Cells.Item(822,8).Value = 8
$ws.Cells.Item(822,9).Formula = "=C822*D822"

# Row 823
$ws.Range("A813:H813").Copy() | Out-Null
$ws.Range("A823:H823").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(823,1).Value = 46063
$ws.Cells.Item(823,2).Value = "Mehdi Boussaid"
$ws.Cells.Item(823,3).Value = 70
$ws.Cells.Item(823,4).Value = 7
$ws.Cells.Item(823,5).Value = 6
$ws.Cells.Item(823,6).Value = 0
$ws.Cells.Item(823,8).Value = 8
$ws.Cells.Item(823,9).Formula = "=C823*D823"

# Row 824
$ws.Range("A813:H813").Copy() | Out-Null
$ws.Range("A824:H824").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(824,1).Value = 46063
$ws.Cells.Item(824,2).Value = "Malik Boussaid"
$ws.Cells.Item(824,3).Value = 70
$ws.Cells.Item(824,4).Value = 4
$ws.Cells.Item(824,5).Value = 2
$ws.Cells.Item(824,6).Value = 0
$ws.Cells.Item(824,8).Value = 10
$ws.Cells.Item(824,9).Formula = "=C824*D824"

# Row 825
$ws.Range("A814:H814").Copy() | Out-Null
$ws.Range("A825:H825").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(825,1).Value = 46063
$ws.Cells.Item(825,2).Value = "Romain Thunet"
$ws.Cells.Item(825,3).Value = 70
$ws.Cells.Item(825,4).Value = 8
$ws.Cells.Item(825,5).Value = 3
$ws.Cells.Item(825,6).Value = 3
$ws.Cells.Item(825,7).Value = "Mollet coup"
$ws.Cells.Item(825,8).Value = 5
$ws.Cells.Item(825,9).Formula = "=C825*D825"

# Row 826
$ws.Range("A814:H814").Copy() | Out-Null
$ws.Range("A826:H826").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(826,1).Value = 46063
$ws.Cells.Item(826,2).Value = "Karahali Souaré"
$ws.Cells.Item(826,3).Value = 70
$ws.Cells.Item(826,4).Value = 9
$ws.Cells.Item(826,5).Value = 9
$ws.Cells.Item(826,6).Value = 6
$ws.Cells.Item(826,7).Value = "Cheville"
$ws.Cells.Item(826,8).Value = 0
$ws.Cells.Item(826,9).Formula = "=C826*D826"

# Row 827
$ws.Range("A813:H813").Copy() | Out-Null
$ws.Range("A827:H827").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(827,1).Value = 46063
$ws.Cells.Item(827,2).Value = "Ilan Ihaddadene"
$ws.Cells.Item(827,3).Value = 70
$ws.Cells.Item(827,4).Value = 8
$ws.Cells.Item(827,5).Value = 8
$ws.Cells.Item(827,6).Value = 0
$ws.Cells.Item(827,8).Value = 6
$ws.Cells.Item(827,9).Formula = "=C827*D827"

# Row 828
$ws.Range("A814:H814").Copy() | Out-Null
$ws.Range("A828:H828").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(828,1).Value = 46063
$ws.Cells.Item(828,2).Value = "Sofiane Belle"
$ws.Cells.Item(828,3).Value = 70
$ws.Cells.Item(828,4).Value = 7
$ws.Cells.Item(828,5).Value = 5
$ws.Cells.Item(828,6).Value = 8
$ws.Cells.Item(828,7).Value = "Côté gauche"
$ws.Cells.Item(828,8).Value = 7
$ws.Cells.Item(828,9).Formula = "=C828*D828"

# Row 829
$ws.Range("A814:H814").Copy() | Out-Null
$ws.Range("A829:H829").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(829,1).Value = 46063
$ws.Cells.Item(829,2).Value = "Hedi Nasri"
$ws.Cells.Item(829,3).Value = 70
$ws.Cells.Item(829,4).Value = 8
$ws.Cells.Item(829,5).Value = 6
$ws.Cells.Item(829,6).Value = 4
$ws.Cells.Item(829,7).Value = "hanche"
$ws.Cells.Item(829,8).Value = 5
$ws.Cells.Item(829,9).Formula = "=C829*D829"

$excel.CutCopyMode = 0

# Restore the scroll position / active selection recorded at save time.
$ws.Activate()
$ws.Range("K825").Select()
$excel.ActiveWindow.ScrollRow = 800
$excel.ActiveWindow.ScrollColumn = 1

